# The workbook contains footnote markers like " [1]".." [5]" appended to several
# vaccine/category names, and a handful of cells use an embedded line break
# (Alt+Enter) to wrap text onto a second line. This edit:
#   1) strips the bracketed footnote markers ("[1]" .. "[5]") from cell text, and
#   2) collapses any embedded line breaks within a cell into a single space,
# across every worksheet in the workbook. Two cells' resulting text becomes an
# exact duplicate of another cell's text elsewhere in the workbook (e.g.
# "Recombivax\nHB" -> "Recombivax HB", which already existed verbatim in another
# cell, and "FluLaval\nQuadrivalent" -> "FluLaval Quadrivalent", likewise) - that
# is expected and is handled automatically by the shared-string table.

$wb = $excel.ActiveWorkbook
$newline = [char]10

foreach ($ws in $wb.Worksheets) {
    $ur = $ws.UsedRange

    # Remove footnote markers "[1]" through "[5]"
    for ($n = 1; $n -le 5; $n++) {
        [void]$ur.Replace("[$n]", "")
    }

    # Collapse embedded line breaks into a single space
    [void]$ur.Replace($newline, " ")
}
